$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datatypes")
$ws.Range("C4").Value = "test"
